$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 3: second term for lease HAN-POISON-001 / 4060-1001-02 ---
# Duplicate row 2 (same lease) into a new row 3, shifting the existing rows down.
$ws.Range("A2:J2").Copy()
$ws.Rows("3:3").Insert()

# Update the new row 3 with the actual term data (start/end dates previous year
# are blank, value previous year missing, new start/end date + value + year).
$ws.Range("C3").ClearContents()
$ws.Range("C3").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Clear()
$ws.Range("F3").Value = 41275
$ws.Range("G3").Value = 41639
$ws.Range("H3").Value = 23000
$ws.Range("I3").Clear()
$ws.Range("J3").Value = 2013

# --- Append a new, empty (formatted only) row at the bottom: row 7 ---
$ws.Range("A2:J2").Copy()
$ws.Range("A7:J7").PasteSpecial(-4122)
$ws.Range("A7:J7").ClearContents()

$excel.CutCopyMode = 0

$ws.Range("L9").Select()
